$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")

# ---------------------------------------------------------------------------
# 1. Fix formatting on row 36 (P36 and AC36 pick up wrap-text formatting to
#    match the rest of the "even" banded rows, e.g. row 34).
# ---------------------------------------------------------------------------
$ws.Range("P36").WrapText = $true
$ws.Range("AC36").WrapText = $true

# ---------------------------------------------------------------------------
# 2. Add new row 37 with the same banded ("odd row") formatting as row 35,
#    then relax the wrap-text columns (M/O/P/AC) to match the no-wrap look
#    used for this new record.
# ---------------------------------------------------------------------------
$ws.Range("A35:AK35").Copy()
$ws.Range("A37:AK37").PasteSpecial(-4122)

$ws.Range("M35").Copy()
$ws.Range("M37,O37,P37,AC37").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Populate the values for the new row (item #35 in the report).
# ---------------------------------------------------------------------------
$ws.Range("A37").Value = 35
$ws.Range("B37").Value = "服務"
$ws.Range("C37").Value = 2025060867
$ws.Range("D37").Value = ""
$ws.Range("E37").Value = ""
$ws.Range("F37").Value = 4125
$ws.Range("G37").Value = "蘆洲長樂店"
$ws.Range("H37").Value = "新北市蘆洲區"
$ws.Range("I37").Value = ""
$ws.Range("J37").Value = ""
$ws.Range("K37").Value = ""
$ws.Range("L37").Value = ""
$ws.Range("M37").Value = ""
$ws.Range("N37").Value = ""
$ws.Range("O37").Value = ""
$ws.Range("P37").Value = ""
$ws.Range("Q37").Value = "THILF04125"
$ws.Range("R37").Value = "新北一"
$ws.Range("S37").Value = "吳宗鴻"
$ws.Range("T37").Value = 1
$ws.Range("U37").Value = "已完工"
$ws.Range("V37").Value = "2025-06-05 17:00:54"
$ws.Range("W37").Value = "2025-06-05 16:10:00"
$ws.Range("X37").Value = "2025-06-05 17:00:00"
$ws.Range("Y37").Value = ""
$ws.Range("Z37").Value = 0.8
$ws.Range("AA37").Value = ""
$ws.Range("AB37").Value = "到場處理"
$ws.Range("AC37").Value = "PMQ2+EDC+STAR"
$ws.Range("AD37").Value = "O"
$ws.Range("AE37").Value = ""
$ws.Range("AF37").Value = ""
$ws.Range("AG37").Value = ""
$ws.Range("AH37").Value = ""
$ws.Range("AI37").Value = ""
$ws.Range("AJ37").Value = "O"
$ws.Range("AK37").Value = "O"

# ---------------------------------------------------------------------------
# 4. Expand the print area to include the new row and move the active
#    selection to the first cell of the new row (mirrors the recorded
#    workbook state after the edit).
# ---------------------------------------------------------------------------
$ws.PageSetup.PrintArea = "$A$1:$AK$37"

$ws.Range("A37").Select()
